# Update the "May" worksheet: renumber the ID column and recompute the
# Amount column for rows 7-61 (the 2nd through 12th "Marty/Mai/Shristi/
# Seble/Rediet" blocks), then adjust the sheet's view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May")

# Column A (ID) becomes a running count (row-1) instead of resetting to
# 1-5 in every 5-row block; Column E (Amount) becomes a straight line
# 34200, 38800, 42400, ... increasing by 4600 each row.
for ($row = 7; $row -le 61; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 5).Value = 34200 + ($row - 7) * 4600
}

# Update the view: zoom to 85% and select E2:E61 (active cell E2).
$excel.ActiveWindow.Zoom = 85
$ws.Range("E2:E61").Select()
